# Motor Templates operating hours update
#
# Rewrites the "(${HR} hours per day, ${DY} days per week, ${WK} weeks per
# year)" sentence to the abbreviated "(${HR} hrs/day, ${DY}  days/wk, ${WK}
# wks/yr)" form, and switches that paragraph's indentation from a
# first-line-only indent to a hanging indent (left=720, hanging=720).

$d = $word.ActiveDocument

# Locate the paragraph that carries the "${HR} ... ${DY} ... ${WK} ..."
# operating-hours sentence (avoid literal "${...}" in a double-quoted
# PowerShell string, which would trigger variable interpolation).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -like '*${HR}*') -and ($t -like '*${DY}*') -and ($t -like '*${WK}*')) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the operating-hours paragraph (containing `${HR}/`${DY}/`${WK})."
}

# Full replacement XML for the paragraph: keeps the original paragraph
# identity/run metadata for the untouched portions, flips <w:ind> to a
# hanging indent, and re-splits the three "per day / per week / per year"
# runs into the abbreviated hrs//day, days//wk, wks//yr runs (with the
# matching <w:proofErr> spell-check bookmarks) called for by the diff.
$targetParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="29955E6B" w14:textId="1EBE0D4F" w:rsidR="001303C6" w:rsidRPr="00646763" w:rsidRDefault="00242BD8" w:rsidP="00242BD8"><w:pPr><w:pStyle w:val="BodyTextIndent"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720" w:hanging="720"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r><w:r w:rsidR="001303C6" w:rsidRPr="00646763"><w:t xml:space="preserve">OH </w:t></w:r><w:r w:rsidR="001303C6" w:rsidRPr="00646763"><w:tab/><w:t xml:space="preserve">= Operating hours of existing fans, </w:t></w:r><w:r w:rsidR="00895696"><w:t>${OH}</w:t></w:r><w:r w:rsidR="001303C6"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001303C6" w:rsidRPr="00646763"><w:t>hrs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001303C6" w:rsidRPr="00646763"><w:t>/yr.</w:t></w:r><w:r w:rsidR="001303C6"><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="2" w:name="_Hlk107346325"/><w:r w:rsidR="001303C6" w:rsidRPr="00833EE1"><w:t>(</w:t></w:r><w:r w:rsidR="00895696"><w:t>${HR}</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>h</w:t></w:r><w:r><w:t>rs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:r><w:t xml:space="preserve">day, </w:t></w:r><w:r w:rsidR="00895696"><w:t>${DY}</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">   </w:t></w:r><w:r><w:t>days</w:t></w:r><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00895696"><w:t>${WK}</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001303C6" w:rsidRPr="00833EE1"><w:t>)</w:t></w:r><w:bookmarkEnd w:id="2"/></w:p>'

$target.Range.InsertXML($targetParaXml)

Write-Host "Operating hours paragraph updated."
